$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 18 data (2027 placeholder predictions; 0s used where data wasn't
# found yet)
$ws.Range("A18").Value = 2027
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 0
$ws.Range("E18").Value = 28653
$ws.Range("F18").Value = 130232
$ws.Range("G18").Value = 4.6
$ws.Range("H18").Value = 28.3
$ws.Range("I18").Value = 30
$ws.Range("J18").Value = 2600

# Extend the D column's ratio formula down through row 18 (mirrors dragging
# the fill handle from D17 to D18), then overwrite D18 with a literal 0
# since B18/C18 are 0 placeholders (would otherwise divide by zero).
$ws.Range("D3:D18").Formula = "=B3/C3"
$ws.Range("D18").Value = 0

$ws.Range("E18").Select() | Out-Null
